$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 679 (pushes existing row 679 and below down by one)
$ws.Rows.Item(679).Insert()

# Populate the newly inserted row 679 with the new log entry
$ws.Cells.Item(679, 1).NumberFormat = "@"
$ws.Cells.Item(679, 1).Value = "2026/01/19"
$ws.Cells.Item(679, 2).Value = "月"
$ws.Cells.Item(679, 3).Value = 4
$ws.Cells.Item(679, 4).Value = 170
